$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (number format/style) of the existing data row (row 2)
# down onto the new row 3 before filling in its values.
$ws.Range("A2:M2").Copy()
$ws.Range("A3:M3").PasteSpecial(-4122)  # xlPasteFormats

# New data row appended to the sheet
$ws.Range("A3").Value = 42600.881111111114
$ws.Range("B3").Value = "Noun"
$ws.Range("C3").Value = 8296
$ws.Range("D3").Value = 8366
$ws.Range("E3").Value = 1465
$ws.Range("F3").Value = 168
$ws.Range("G3").Value = 110
$ws.Range("H3").Value = 58
$ws.Range("I3").Value = 38
$ws.Range("J3").Value = 2
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 100
$ws.Range("M3").Value = 0

# Widen column A slightly, as in the diff (13.85546875 -> 14.85546875)
$ws.Columns.Item(1).ColumnWidth = 14
